# Sprint 4 Backlog - Burndown: record actual times for remaining tasks
# (rows 15, 17, 18, 19, 20 on Sheet1) per the commit "Updated the sprint
# backlog to include my times / All my times are now updated".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 - "Add new recipe / Complete functionality add an new recipe (web)"
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = "Matthew"
$ws.Range("I15").Value = 0

# Row 17 - "Share a recipe / Complete functionality to share a recipe (web)"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = "Matthew"
$ws.Range("H17").Value = 0.5
$ws.Range("I17").Value = 0

# Row 18 - "Share a recipe / Complete functionality to share a recipe (desktop)"
$ws.Range("C18").Value = 1.5
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = "Destiny"
$ws.Range("H18").Value = 1.5
$ws.Range("I18").Value = 0

# Row 19 - "Share a recipe / Create UI to view shared recipes (desktop)"
$ws.Range("E19").Value = 0.5
$ws.Range("F19").Value = "Matthew"
$ws.Range("H19").Value = 0.5
$ws.Range("I19").Value = 0

# Row 20 - "Share a recipe / Create UI to view shared recipes (web)"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "Destiny"
$ws.Range("I20").Value = 0

# Recalculate totals / SUMIFs / the embedded burndown chart's source cells
$excel.Calculate()

$wb.Save()
